# Generate Report for Handoff
# Adds two newly-handed-off files (93e4a31f... and 99ead3f0...) to the
# localization status report, pushing the ".localization-config" row down
# on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper data
# ---------------------------------------------------------------------
$guid1 = "93e4a31f-1314-42f3-8e56-bb922eefd6e8"
$guid2 = "99ead3f0-7917-49df-9654-a4f94c45a773"

$guid1Md  = "$guid1.md"
$guid2Md  = "$guid2.md"

$guid1ZhCn = "$guid1.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.zh-cn.xlf"
$guid1DeDe = "$guid1.d37e440889d76a4c11dbd6d0fabd0f9c5cd90db5.de-de.xlf"
$guid2ZhCn = "$guid2.7b5ccef2be4187d8d14e2b451a398be9125d803f.zh-cn.xlf"
$guid2DeDe = "$guid2.7b5ccef2be4187d8d14e2b451a398be9125d803f.de-de.xlf"

$guid1MdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/d3850476015e46fbd073fbcc4937803b51b3335e/e2e/$guid1Md"
$guid2MdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/d3850476015e46fbd073fbcc4937803b51b3335e/e2e/$guid2Md"
$configUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/d3850476015e46fbd073fbcc4937803b51b3335e/.localization-config"

$guid1ZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/932cf51e579bd44d8c2f348d30702cd3462fda5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$guid1ZhCn"
$guid2ZhCnUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/932cf51e579bd44d8c2f348d30702cd3462fda5e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$guid2ZhCn"
$guid1DeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02a6f8a826f188add1b3d1c47d6bfc154cb262e9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$guid1DeDe"
$guid2DeDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02a6f8a826f188add1b3d1c47d6bfc154cb262e9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$guid2DeDe"

$readyStatus = "Ready for handoff"
$notLocalized = "Not to be localized"
$configName = ".localization-config"

$handoffDateZhCn = "2016-02-25 06:45:03"
$handoffDateDeDe = "2016-02-25 06:45:17"
$zeroDate = "0001-01-01 00:00:00"
$includeStatus = "Include"
$ignoredStatus = "Ignored"

$hyperlinkFontColor = 15570276

function Set-HyperlinkFont($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkFontColor
}

function Replace-Hyperlink($ws, $cellAddr, $displayText, $url) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.Delete()
        }
    }
    $ws.Range($cellAddr.Replace('$', '')).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellAddr.Replace('$', '')), $url, "", "", $displayText)
    Set-HyperlinkFont($ws.Range($cellAddr.Replace('$', '')))
}

function Add-Hyperlink($ws, $cellRef, $displayText, $url) {
    $ws.Range($cellRef).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $displayText)
    Set-HyperlinkFont($ws.Range($cellRef))
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 4 becomes the first new file (was ".localization-config")
Replace-Hyperlink $ws1 '$A$4' $guid1Md $guid1MdUrl
$ws1.Range("B4").Value = $readyStatus
$ws1.Range("C4").Value = $readyStatus

# Row 5: second new file
Add-Hyperlink $ws1 "A5" $guid2Md $guid2MdUrl
$ws1.Range("B5").Value = $readyStatus
$ws1.Range("C5").Value = $readyStatus

# Row 6: ".localization-config" moves down here
Add-Hyperlink $ws1 "A6" $configName $configUrl
$ws1.Range("B6").Value = $notLocalized
$ws1.Range("C6").Value = $notLocalized

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 4: fill in handoff file + datetime for the first new file,
# and point A4 at the new file instead of the config file.
Replace-Hyperlink $ws2 '$A$4' $guid1Md $guid1MdUrl
$ws2.Range("B4").Value = $readyStatus
Add-Hyperlink $ws2 "C4" $guid1ZhCn $guid1ZhCnUrl
$ws2.Range("D4").Value = $handoffDateZhCn
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = $zeroDate
$ws2.Range("H4").Value = $includeStatus

# Row 5: second new file
Add-Hyperlink $ws2 "A5" $guid2Md $guid2MdUrl
$ws2.Range("B5").Value = $readyStatus
Add-Hyperlink $ws2 "C5" $guid2ZhCn $guid2ZhCnUrl
$ws2.Range("D5").Value = $handoffDateZhCn
$ws2.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G5").Value = $zeroDate
$ws2.Range("H5").Value = $includeStatus

# Row 6: ".localization-config" moves down here
Add-Hyperlink $ws2 "A6" $configName $configUrl
$ws2.Range("B6").Value = $notLocalized
$ws2.Range("D6").Value = $zeroDate
$ws2.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G6").Value = $zeroDate
$ws2.Range("H6").Value = $ignoredStatus

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

Replace-Hyperlink $ws3 '$A$4' $guid1Md $guid1MdUrl
$ws3.Range("B4").Value = $readyStatus
Add-Hyperlink $ws3 "C4" $guid1DeDe $guid1DeDeUrl
$ws3.Range("D4").Value = $handoffDateDeDe
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = $zeroDate
$ws3.Range("H4").Value = $includeStatus

Add-Hyperlink $ws3 "A5" $guid2Md $guid2MdUrl
$ws3.Range("B5").Value = $readyStatus
Add-Hyperlink $ws3 "C5" $guid2DeDe $guid2DeDeUrl
$ws3.Range("D5").Value = $handoffDateDeDe
$ws3.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G5").Value = $zeroDate
$ws3.Range("H5").Value = $includeStatus

Add-Hyperlink $ws3 "A6" $configName $configUrl
$ws3.Range("B6").Value = $notLocalized
$ws3.Range("D6").Value = $zeroDate
$ws3.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G6").Value = $zeroDate
$ws3.Range("H6").Value = $ignoredStatus

Write-Host "Report generated for handoff."
